# Add columns I (I0) and J (IF) to the sheet, mirroring the style of the
# existing header (H1 "IP") and filling data rows 2-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: copy formatting from the existing header cell H1 so the new
# headers reuse the same cell style (bold/border/centered) as the rest of
# row 1, instead of Excel minting a brand-new style entry.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-39 for columns I and J.
$data = @{
    2  = @(9, 9)
    3  = @(6, 6)
    4  = @(8, 8)
    5  = @(8, 9)
    6  = @(6, 6)
    7  = @(7, 7)
    8  = @(4, 5)
    9  = @(6, 6)
    10 = @(7, 7)
    11 = @(5, 6)
    12 = @(5, 5)
    13 = @(6, 6)
    14 = @(3, 4)
    15 = @(5, 5)
    16 = @(7, 7)
    17 = @(5, 5)
    18 = @(8, 8)
    19 = @(9, 9)
    20 = @(5, 5)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(7, 8)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(7, 7)
    29 = @(6, 6)
    30 = @(7, 7)
    31 = @(6, 7)
    32 = @(7, 7)
    33 = @(6, 6)
    34 = @(5, 5)
    35 = @(1, 2)
    36 = @(5, 6)
    37 = @(6, 6)
    38 = @(5, 5)
    39 = @(3, 3)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
